$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell address -> new numeric value, derived from the daily TGP price-sheet refresh.
$updates = @{
    "A8" = 45967;
    "D8" = 168.28;
    "E8" = 160.72999999999999;
    "F8" = 170.73;
    "G8" = 160.9;
    "A9" = 45967;
    "D9" = 168.28;
    "E9" = 160.72999999999999;
    "F9" = 170.73;
    "G9" = 160.9;
    "A10" = 45967;
    "D10" = 170.11;
    "E10" = 162.91999999999999;
    "F10" = 172.92;
    "G10" = 163.41;
    "A11" = 45966;
    "D11" = 168.21;
    "E11" = 160.83000000000001;
    "F11" = 170.83;
    "G11" = 161;
    "A12" = 45966;
    "D12" = 168.21;
    "E12" = 160.83000000000001;
    "F12" = 170.83;
    "G12" = 161;
    "A13" = 45966;
    "D13" = 170.13;
    "E13" = 163.01;
    "F13" = 173.01;
    "G13" = 163.51;
    "A17" = 45967;
    "D17" = 173.61;
    "E17" = 165.89;
    "F17" = 175.89;
    "A18" = 45966;
    "D18" = 173.64;
    "E18" = 165.97;
    "F18" = 175.97;
    "A22" = 45967;
    "D22" = 169.31;
    "E22" = 161.94999999999999;
    "F22" = 171.55;
    "G22" = 163.22999999999999;
    "A23" = 45967;
    "D23" = 174.88;
    "E23" = 166.65;
    "F23" = 176.65;
    "A24" = 45967;
    "D24" = 174.69;
    "E24" = 166.82;
    "F24" = 176.82;
    "A25" = 45967;
    "D25" = 175.52;
    "E25" = 166.21;
    "F25" = 176.21;
    "G25" = 166.25;
    "A26" = 45967;
    "D26" = 174.26;
    "E26" = 167.77;
    "F26" = 177.77;
    "A27" = 45966;
    "D27" = 169.34;
    "E27" = 162.04;
    "F27" = 171.64;
    "G27" = 163.33000000000001;
    "A28" = 45966;
    "D28" = 174.9;
    "E28" = 166.74;
    "F28" = 176.74;
    "A29" = 45966;
    "D29" = 174.72;
    "E29" = 166.91;
    "F29" = 176.91;
    "A30" = 45966;
    "D30" = 175.55;
    "E30" = 166.3;
    "F30" = 176.3;
    "G30" = 166.34;
    "A31" = 45966;
    "D31" = 174.29;
    "E31" = 167.86;
    "F31" = 177.86;
    "A35" = 45967;
    "D35" = 168.69;
    "E35" = 160.13999999999999;
    "F35" = 169.14;
    "A36" = 45966;
    "D36" = 168.61;
    "E36" = 160.22999999999999;
    "F36" = 169.23;
    "A40" = 45967;
    "D40" = 173.96;
    "E40" = 165.52;
    "F40" = 175.52;
    "A41" = 45967;
    "D41" = 173.67;
    "E41" = 165.94;
    "F41" = 175.95;
    "A42" = 45966;
    "D42" = 173.98;
    "E42" = 165.59;
    "F42" = 175.59;
    "A43" = 45966;
    "D43" = 173.69;
    "E43" = 166.01;
    "F43" = 176.01;
    "A47" = 45967;
    "D47" = 168.34;
    "E47" = 161.63999999999999;
    "F47" = 171.64;
    "A48" = 45967;
    "D48" = 168.36;
    "E48" = 161.83000000000001;
    "F48" = 171.83;
    "A49" = 45966;
    "D49" = 168.62;
    "E49" = 161.96;
    "F49" = 171.96;
    "A50" = 45966;
    "D50" = 168.63;
    "E50" = 162.15;
    "F50" = 172.15;
    "A54" = 45967;
    "D54" = 184.18;
    "E54" = 175.67;
    "F54" = 185.67;
    "A55" = 45967;
    "D55" = 171.85;
    "E55" = 173.59;
    "F55" = 183.59;
    "A56" = 45967;
    "D56" = 174.33;
    "A57" = 45967;
    "D57" = 174.09;
    "E57" = 167.85;
    "A58" = 45967;
    "D58" = 170;
    "E58" = 163.91;
    "F58" = 173.91;
    "A59" = 45967;
    "D59" = 176.61;
    "E59" = 173.95;
    "A60" = 45966;
    "D60" = 184.2;
    "A61" = 45966;
    "D61" = 171.87;
    "E61" = 173.58;
    "F61" = 183.58;
    "A62" = 45966;
    "D62" = 174.36;
    "A63" = 45966;
    "D63" = 174.12;
    "E63" = 167.84;
    "A64" = 45966;
    "D64" = 170.03;
    "E64" = 163.9;
    "F64" = 173.9;
    "A65" = 45966;
    "D65" = 176.65;
    "E65" = 174.03
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
